$d = $word.ActiveDocument

# 1. "Super admin." -> "//Super admin."
$d.Content.Find.Execute("Super admin.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "//Super admin.", 2)

# 2. "Admin." -> "//Admin."
$d.Content.Find.Execute("Admin.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "//Admin.", 2)

# 3. Insert a new paragraph with text after the "//Pictures sirf..." paragraph
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Pictures sirf and more details baaki all others in details.php*") {
        $target = $p
        break
    }
}

$newPara = $target.Range.InsertParagraphAfter()
$newRange = $target.Next().Range
$newRange.Text = "Alhamdulliah i am satisfied with what i have so far and i've worked really hard on this."
